$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.305838942527771
$ws.Range("B1").Value = 2.183547496795654
$ws.Range("C1").Value = 4.857301712036133
$ws.Range("D1").Value = 3.165481567382812
$ws.Range("E1").Value = 1.331758975982666
